$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 175
$ws.Range('B175').Value = 7302795
$ws.Range('E175').Value = 'Unin Comercio'
$ws.Range('F175').Value = 'Deportivo Garcilaso'
$ws.Range('G175').Value = 1
$ws.Range('H175').Value = 2
$ws.Range('I175').Value = 1
$ws.Range('J175').Value = 0
$ws.Range('K175').Value = 'A'
$ws.Range('L175').Value = 2.25
$ws.Range('M175').Value = 3.3
$ws.Range('N175').Value = 2.7
$ws.Range('O175').Value = 1.75
$ws.Range('P175').Value = 3.6
$ws.Range('Q175').Value = 4
$ws.Range('R175').Value = -0.5
$ws.Range('S175').Value = 1.8
$ws.Range('T175').Value = 2
$ws.Range('U175').Value = 2.75
$ws.Range('V175').Value = 1.825
$ws.Range('W175').Value = 1.975
$ws.Range('X175').Value = -1
$ws.Range('Y175').Value = -1
$ws.Range('Z175').Value = 3
$ws.Range('AA175').Value = -1
$ws.Range('AB175').Value = 1
$ws.Range('AC175').Value = 0.4125
$ws.Range('AD175').Value = -0.5

# Row 176
$ws.Range('B176').Value = 7302200
$ws.Range('E176').Value = 'Carlos Manucci'
$ws.Range('F176').Value = 'Deportivo Binacional'
$ws.Range('G176').Value = 3
$ws.Range('H176').Value = 2
$ws.Range('I176').Value = 0
$ws.Range('J176').Value = 0
$ws.Range('K176').Value = 'H'
$ws.Range('L176').Value = 2
$ws.Range('M176').Value = 3.2
$ws.Range('N176').Value = 3.75
$ws.Range('O176').Value = 1.75
$ws.Range('P176').Value = 3.4
$ws.Range('Q176').Value = 4.333
$ws.Range('R176').Value = -0.5
$ws.Range('S176').Value = 1.85
$ws.Range('T176').Value = 1.95
$ws.Range('U176').Value = 2.5
$ws.Range('V176').Value = 1.85
$ws.Range('W176').Value = 1.95
$ws.Range('X176').Value = 0.75
$ws.Range('Y176').Value = -1
$ws.Range('Z176').Value = -1
$ws.Range('AA176').Value = 0.8500000000000001
$ws.Range('AB176').Value = -1
$ws.Range('AC176').Value = 0.8500000000000001
$ws.Range('AD176').Value = -1

# Row 184
$ws.Range('B184').Value = 7384626
$ws.Range('E184').Value = 'Sporting Cristal'
$ws.Range('F184').Value = 'Alianza Atletico'
$ws.Range('G184').Value = 3
$ws.Range('H184').Value = 0
$ws.Range('I184').Value = 3
$ws.Range('J184').Value = 0
$ws.Range('K184').Value = 'H'
$ws.Range('L184').Value = 1.3
$ws.Range('M184').Value = 5
$ws.Range('N184').Value = 9
$ws.Range('O184').Value = 1.166
$ws.Range('P184').Value = 6.5
$ws.Range('Q184').Value = 13
$ws.Range('R184').Value = -2
$ws.Range('S184').Value = 1.85
$ws.Range('T184').Value = 1.95
$ws.Range('U184').Value = 3.25
$ws.Range('V184').Value = 2
$ws.Range('W184').Value = 1.8
$ws.Range('X184').Value = 0.1659999999999999
$ws.Range('Y184').Value = -1
$ws.Range('Z184').Value = -1
$ws.Range('AA184').Value = 0.8500000000000001
$ws.Range('AB184').Value = -1
$ws.Range('AC184').Value = -0.5
$ws.Range('AD184').Value = 0.4

# Row 185
$ws.Range('B185').Value = 7384627
$ws.Range('E185').Value = 'Universitario de Deportes'
$ws.Range('F185').Value = 'Sport Huancayo'
$ws.Range('G185').Value = 2
$ws.Range('H185').Value = 0
$ws.Range('I185').Value = 1
$ws.Range('J185').Value = 0
$ws.Range('K185').Value = 'H'
$ws.Range('L185').Value = 1.25
$ws.Range('M185').Value = 5
$ws.Range('N185').Value = 12
$ws.Range('O185').Value = 1.181
$ws.Range('P185').Value = 6
$ws.Range('Q185').Value = 13
$ws.Range('R185').Value = -1.75
$ws.Range('S185').Value = 1.8
$ws.Range('T185').Value = 2
$ws.Range('U185').Value = 2.75
$ws.Range('V185').Value = 1.85
$ws.Range('W185').Value = 1.95
$ws.Range('X185').Value = 0.181
$ws.Range('Y185').Value = -1
$ws.Range('Z185').Value = -1
$ws.Range('AA185').Value = 0.4
$ws.Range('AB185').Value = -0.5
$ws.Range('AC185').Value = -1
$ws.Range('AD185').Value = 0.95

# Row 186
$ws.Range('B186').Value = 7384629
$ws.Range('E186').Value = 'Deportivo Garcilaso'
$ws.Range('F186').Value = 'Alianza Lima'
$ws.Range('G186').Value = 0
$ws.Range('H186').Value = 1
$ws.Range('I186').Value = 0
$ws.Range('J186').Value = 1
$ws.Range('K186').Value = 'A'
$ws.Range('L186').Value = 2.625
$ws.Range('M186').Value = 3.3
$ws.Range('N186').Value = 2.5
$ws.Range('O186').Value = 2.7
$ws.Range('P186').Value = 3.4
$ws.Range('Q186').Value = 2.375
$ws.Range('R186').Value = 0
$ws.Range('S186').Value = 2.025
$ws.Range('T186').Value = 1.775
$ws.Range('U186').Value = 2.25
$ws.Range('V186').Value = 1.825
$ws.Range('W186').Value = 1.975
$ws.Range('X186').Value = -1
$ws.Range('Y186').Value = -1
$ws.Range('Z186').Value = 1.375
$ws.Range('AA186').Value = -1
$ws.Range('AB186').Value = 0.7749999999999999
$ws.Range('AC186').Value = -1
$ws.Range('AD186').Value = 0.9750000000000001

# Row 187
$ws.Range('B187').Value = 7384628
$ws.Range('E187').Value = 'Deportivo Binacional'
$ws.Range('F187').Value = 'FBC Melgar'
$ws.Range('G187').Value = 1
$ws.Range('H187').Value = 2
$ws.Range('I187').Value = 1
$ws.Range('J187').Value = 1
$ws.Range('K187').Value = 'A'
$ws.Range('L187').Value = 2.75
$ws.Range('M187').Value = 3.3
$ws.Range('N187').Value = 2.375
$ws.Range('O187').Value = 3.3
$ws.Range('P187').Value = 3.6
$ws.Range('Q187').Value = 2
$ws.Range('R187').Value = 0.5
$ws.Range('S187').Value = 1.8
$ws.Range('T187').Value = 2
$ws.Range('U187').Value = 2.75
$ws.Range('V187').Value = 1.975
$ws.Range('W187').Value = 1.875
$ws.Range('X187').Value = -1
$ws.Range('Y187').Value = -1
$ws.Range('Z187').Value = 1
$ws.Range('AA187').Value = -1
$ws.Range('AB187').Value = 1
$ws.Range('AC187').Value = 0.4875
$ws.Range('AD187').Value = -0.5

# Row 188
$ws.Range('B188').Value = 7384625
$ws.Range('E188').Value = 'AD Tarma'
$ws.Range('F188').Value = 'Carlos Manucci'
$ws.Range('G188').Value = 0
$ws.Range('H188').Value = 0
$ws.Range('I188').Value = 0
$ws.Range('J188').Value = 0
$ws.Range('K188').Value = 'D'
$ws.Range('L188').Value = 1.5
$ws.Range('M188').Value = 3.75
$ws.Range('N188').Value = 7
$ws.Range('O188').Value = 1.363
$ws.Range('P188').Value = 4.333
$ws.Range('Q188').Value = 9.5
$ws.Range('R188').Value = -1.25
$ws.Range('S188').Value = 1.875
$ws.Range('T188').Value = 1.925
$ws.Range('U188').Value = 2.5
$ws.Range('V188').Value = 1.8
$ws.Range('W188').Value = 2
$ws.Range('X188').Value = -1
$ws.Range('Y188').Value = 3.333
$ws.Range('Z188').Value = -1
$ws.Range('AA188').Value = -1
$ws.Range('AB188').Value = 0.925
$ws.Range('AC188').Value = -1
$ws.Range('AD188').Value = 1

# Row 228
$ws.Range('B228').Value = 7818817
$ws.Range('E228').Value = 'Sport Boys'
$ws.Range('F228').Value = 'Cusco FC'
$ws.Range('G228').Value = 3
$ws.Range('H228').Value = 0
$ws.Range('I228').Value = 2
$ws.Range('J228').Value = 0
$ws.Range('K228').Value = 'H'
$ws.Range('L228').Value = 2.2
$ws.Range('M228').Value = 3.2
$ws.Range('N228').Value = 3.2
$ws.Range('O228').Value = 1.6
$ws.Range('P228').Value = 3.75
$ws.Range('Q228').Value = 5.75
$ws.Range('R228').Value = -0.75
$ws.Range('S228').Value = 1.85
$ws.Range('T228').Value = 2
$ws.Range('U228').Value = 2.5
$ws.Range('V228').Value = 1.975
$ws.Range('W228').Value = 1.875
$ws.Range('X228').Value = 0.6000000000000001
$ws.Range('Y228').Value = -1
$ws.Range('Z228').Value = -1
$ws.Range('AA228').Value = 0.8500000000000001
$ws.Range('AB228').Value = -1
$ws.Range('AC228').Value = 0.9750000000000001
$ws.Range('AD228').Value = -1

# Row 229
$ws.Range('B229').Value = 7818816
$ws.Range('E229').Value = 'UTC Cajamarca'
$ws.Range('F229').Value = 'Universitario de Deportes'
$ws.Range('G229').Value = 0
$ws.Range('H229').Value = 0
$ws.Range('I229').Value = 0
$ws.Range('J229').Value = 0
$ws.Range('K229').Value = 'D'
$ws.Range('L229').Value = 3.3
$ws.Range('M229').Value = 3.3
$ws.Range('N229').Value = 2.1
$ws.Range('O229').Value = 4.5
$ws.Range('P229').Value = 3.2
$ws.Range('Q229').Value = 1.95
$ws.Range('R229').Value = 0.5
$ws.Range('S229').Value = 2
$ws.Range('T229').Value = 1.85
$ws.Range('U229').Value = 2
$ws.Range('V229').Value = 1.775
$ws.Range('W229').Value = 2.1
$ws.Range('X229').Value = -1
$ws.Range('Y229').Value = 2.2
$ws.Range('Z229').Value = -1
$ws.Range('AA229').Value = 1
$ws.Range('AB229').Value = -1
$ws.Range('AC229').Value = -1
$ws.Range('AD229').Value = 1.1
